$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered, centered) from H1 to I1:J1, matching existing header formatting
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF) for rows 2-54
$iVals = @(7,8,8,8,7,6,7,9,5,7,7,9,7,8,3,9,6,9,7,4,9,7,9,8,8,9,7,5,7,8,5,7,7,10,5,10,5,8,5,7,9,8,6,7,7,6,8,6,7,5,3,4,3)
$jVals = @(8,8,8,8,9,8,8,9,6,7,8,9,8,8,4,9,6,9,8,6,10,7,9,8,9,9,8,6,7,8,7,7,7,10,6,10,6,8,6,8,9,9,7,8,8,6,9,6,8,6,4,4,3)
for ($r = 2; $r -le 54; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}

$excel.CutCopyMode = $false
